# Append one new daily-portfolio row (row 44) to Sheet1, extending the
# existing A1:D43 table of Date / SUZLON.NS / TATAMOTORS.NS / ETERNAL.NS
# by a single day (2025-09-28), mirroring the pattern of the prior rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 44

# Column A holds the date as literal text (e.g. "2025-09-27" in the row
# above), not a real Excel date serial. Force the cell to Text format
# before assignment so the COM layer doesn't auto-convert the
# date-shaped string into a date serial number, then drop the style
# back to Normal so the cell matches the unstyled cells above it.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-09-28"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 55.27999877929688
$ws.Cells.Item($row, 3).Value = 672.9000244140625
$ws.Cells.Item($row, 4).Value = 321
